$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row of data to the work-log table (row 11)
$ws.Range("A11").Value = "牟秋宇"
$ws.Range("B11").Value = "参与前端代码编写"
$ws.Range("C11").Value = 600
$ws.Range("E11").Value = 6

$ws.Range("F11").Select()
